$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (M2)
$ws.Range("B3").Value = -1234.747413556725
$ws.Range("C3").Value = 2085.918756089357
$ws.Range("D3").Value = 255.3567153003003
$ws.Range("F3").Value = 2531.494827113449
$ws.Range("G3").Value = 2641.273997013361

# Row 4 (M3)
$ws.Range("B4").Value = -1173.327028153149
$ws.Range("C4").Value = 534.3381372048618
$ws.Range("D4").Value = 210.2897315333678
$ws.Range("F4").Value = 2440.654056306297
$ws.Range("G4").Value = 2607.093442928743

# Row 5 (M4)
$ws.Range("B5").Value = -1140.876482816072
$ws.Range("C5").Value = 324.3867877354729
$ws.Range("D5").Value = 183.9783236976492
$ws.Range("F5").Value = 2407.752965632144
$ws.Range("G5").Value = 2630.852568977125
$ws.Range("H5").Value = 0.000000007474498020521025

# Row 6 (M5)
$ws.Range("B6").Value = -1117.80714587823
$ws.Range("C6").Value = 269.2398492068251
$ws.Range("D6").Value = 170.426483028431
$ws.Range("F6").Value = 2393.61429175646
$ws.Range("G6").Value = 2673.374111823976
$ws.Range("H6").Value = 0.000007594868898008755
